# Update from PC 280423 2100 Final Configuration files.
# Insert new "exceptions" rows into the חריגים sheet, interleaved with the
# pre-existing rows (which shift down but keep their own cell contents).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to force a purely-numeric-looking string to be stored as TEXT
# (shared string) instead of being auto-coerced to a number, while leaving
# the cell's style back at the default (no visible formatting change).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# 1) Insert the six new rows at their final positions, top to bottom, so the
#    three original rows (ברונפלד איתן / גייר גפן / גרדי דן) shift down to
#    rows 3, 6 and 9 while keeping their existing cell values untouched.
$ws.Rows("2:2").Insert()
$ws.Rows("4:4").Insert()
$ws.Rows("5:5").Insert()
$ws.Rows("7:7").Insert()
$ws.Rows("8:8").Insert()
$ws.Rows("10:10").Insert()

# 2) Fill in the data for each newly inserted (blank) row.

# Row 2: בר הלוי יורם
Set-TextValue $ws.Range("A2") "56528482"
$ws.Range("B2").Value = "בר הלוי יורם"
$ws.Range("D2").Value = "קצר"
Set-TextValue $ws.Range("E2") "503147001"
$ws.Range("G2").Value = "נרשם לתחרותי אבל ללא כרטיס אלקטרוני"

# Row 4: גוטהילף צבי
Set-TextValue $ws.Range("A4") "52674090"
$ws.Range("B4").Value = "גוטהילף צבי"
$ws.Range("D4").Value = "בינוני"
Set-TextValue $ws.Range("E4") "543295075"
$ws.Range("G4").Value = "נרשם לתחרותי אבל ללא כרטיס אלקטרוני"

# Row 5: גורקה אורנה
Set-TextValue $ws.Range("A5") "22512727"
$ws.Range("B5").Value = "גורקה אורנה"
$ws.Range("D5").Value = "קצרצר"
Set-TextValue $ws.Range("E5") "543250080"
$ws.Range("G5").Value = "נרשם לתחרותי אבל ללא כרטיס אלקטרוני"

# Row 7: כהן ערן
Set-TextValue $ws.Range("A7") "24219149"
$ws.Range("B7").Value = "כהן ערן"
$ws.Range("D7").Value = "קצר"
Set-TextValue $ws.Range("E7") "546644385"
$ws.Range("G7").Value = "נרשם לתחרותי אבל ללא כרטיס אלקטרוני"

# Row 8: בן זאב אורי
Set-TextValue $ws.Range("A8") "303054878"
$ws.Range("B8").Value = "בן זאב אורי"
$ws.Range("D8").Value = "בינוני"
Set-TextValue $ws.Range("E8") "526868697"
Set-TextValue $ws.Range("F8") "7010070"
$ws.Range("G8").Value = "רשם הערה בעת ההרשמה לתחרות"

# Row 10: צופה אלעד
Set-TextValue $ws.Range("A10") "6744"
$ws.Range("B10").Value = "צופה אלעד"
$ws.Range("C10").Value = "טכניון כרמל"
$ws.Range("D10").Value = "קצרצר"
$ws.Range("E10").Value = "052-4044409"
Set-TextValue $ws.Range("F10") "2071510"
$ws.Range("G10").Value = "רשם הערה בעת ההרשמה לתחרות"
